# Fruta / hortaliza, semanal
# Update the weekly data rows (2-14) with the new permutation of values
# for columns D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2  = @{ D = 44251; J = 120; K = 5000; L = 5000; M = 5000; O = "Región Metropolitana"; P = 312 }
  3  = @{ D = 44187; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  4  = @{ D = 44186; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  5  = @{ D = 44232; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  6  = @{ D = 44292; J = 90;  K = 6000; L = 6000; M = 6000; O = "Región Metropolitana"; P = 375 }
  7  = @{ D = 44188; J = 210; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  8  = @{ D = 44208; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  9  = @{ D = 44210; J = 340; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  10 = @{ D = 44231; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  11 = @{ D = 44215; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  12 = @{ D = 44204; J = 430; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  13 = @{ D = 44230; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  14 = @{ D = 44189; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
